$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Ancient Arena 930PM" event row (row 11, column A only)
$ws.Range("A11").Value = "Ancient Arena 930PM"

# Match the new active selection left by the editing session
$ws.Range("D12").Select() | Out-Null
